# "hospital planner closes at night"
# The Hospital sheet's overnight hour slots (previously "livre"/free, plus
# two special multi-booking slots) are marked as "fechado" (closed).
# This also causes the now-unused "pediatria-Patient1;cardiologia-Patient3"
# and "pediatria-Patient1" shared strings to be dropped, and the new
# "fechado" string to be appended, which naturally re-indexes the
# "marcado"/"ocupador" references on the other patient sheets.

$wb = $excel.ActiveWorkbook
$hospital = $wb.Worksheets.Item("Hospital")
$patient2 = $wb.Worksheets.Item("Patient2")

# Overnight / closed blocks on the Hospital timetable.
$hospital.Range("B1:B9").Value = "fechado"
$hospital.Range("B22:B34").Value = "fechado"
$hospital.Range("B47:B59").Value = "fechado"
$hospital.Range("B72:B83").Value = "fechado"
$hospital.Range("B97:B100").Value = "fechado"

# View-state: Hospital becomes the active tab/sheet, with a new selection
# further down the sheet; Patient2 (previously the active tab) loses the
# selection it had and resets to B1.
$patient2.Range("B1").Select()
$hospital.Activate()
$hospital.Range("B103").Select()
